$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row text (D1 first so "Selections (...)" gets the lower
# new shared-string index, then A1 so "PlaceID..." gets the next one) ---
$ws.Range("D1").Value = "Selections ( 빈칸없이 )"
$ws.Range("A1").Value = "PlaceID`n일반맵 <10000`n특별맵 >10000"

# --- Column widths (A is new, B grows a bit to fit the new header text) ---
$ws.Columns.Item(1).ColumnWidth = 15.5
$ws.Columns.Item(2).ColumnWidth = 16.75

# --- Header row formatting: A1/C1 centered + wrapped, B1/D1 centered ---
$xlCenter = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter

foreach ($addr in @("A1", "C1")) {
    $ws.Range($addr).HorizontalAlignment = $xlCenter
    $ws.Range($addr).VerticalAlignment = $xlCenter
    $ws.Range($addr).WrapText = $true
}

foreach ($addr in @("B1", "D1")) {
    $ws.Range($addr).HorizontalAlignment = $xlCenter
    $ws.Range($addr).VerticalAlignment = $xlCenter
    $ws.Range($addr).WrapText = $false
}

# --- Row height for the now taller, wrapped header row ---
$ws.Rows.Item(1).RowHeight = 49.5

# --- Selection moves to C1 ---
$ws.Range("C1").Select() | Out-Null
